$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.016.58'
$ws.Range('E2').Value = '  -2.78%  '
$ws.Range('D3').Value = '1.890.00'
$ws.Range('E3').Value = '  -3.16%  '
$s = $ws.Range('D4').Style
$ws.Range('D4').Value = "'0.9994"
$ws.Range('D4').Style = $s
$ws.Range('E4').Value = '  -0.10%  '
$s = $ws.Range('D5').Style
$ws.Range('D5').Value = "'329.75"
$ws.Range('D5').Style = $s
$ws.Range('E5').Value = '  -3.82%  '
$s = $ws.Range('D6').Style
$ws.Range('D6').Value = "'0.9992"
$ws.Range('D6').Style = $s
$ws.Range('E6').Value = '  -0.04%  '
$s = $ws.Range('D7').Style
$ws.Range('D7').Value = "'0.4590"
$ws.Range('D7').Style = $s
$ws.Range('E7').Value = '  -4.15%  '
$s = $ws.Range('D8').Style
$ws.Range('D8').Value = "'0.4105"
$ws.Range('D8').Style = $s
$ws.Range('E8').Value = '  -1.33%  '
$s = $ws.Range('D9').Style
$ws.Range('D9').Value = "'47.76"
$ws.Range('D9').Style = $s
$ws.Range('E9').Value = '  -2.19%  '
$s = $ws.Range('D10').Style
$ws.Range('D10').Value = "'0.07959"
$ws.Range('D10').Style = $s
$ws.Range('E10').Value = '  -3.87%  '
$s = $ws.Range('D11').Style
$ws.Range('D11').Value = "'0.9980"
$ws.Range('D11').Style = $s
$ws.Range('E11').Value = '  -4.80%  '
$s = $ws.Range('D12').Style
$ws.Range('D12').Value = "'21.74"
$ws.Range('D12').Style = $s
$ws.Range('E12').Value = '  -4.57%  '
$ws.Range('D13').Value = '1.864.24'
$ws.Range('E13').Value = '  -4.15%  '
$s = $ws.Range('D14').Style
$ws.Range('D14').Value = "'5.920"
$ws.Range('D14').Style = $s
$ws.Range('E14').Value = '  -4.47%  '
$s = $ws.Range('D15').Style
$ws.Range('D15').Value = "'7.075"
$ws.Range('D15').Style = $s
$ws.Range('E15').Value = '  -5.22%  '
$s = $ws.Range('D16').Style
$ws.Range('D16').Value = "'0.9997"
$ws.Range('D16').Style = $s
$ws.Range('E16').Value = '  -0.16%  '
$s = $ws.Range('D17').Style
$ws.Range('D17').Value = "'88.59"
$ws.Range('D17').Style = $s
$ws.Range('E17').Value = '  -4.71%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$s = $ws.Range('D18').Style
$ws.Range('D18').Value = "'0.00001026"
$ws.Range('D18').Style = $s
$ws.Range('E18').Value = '  -3.98%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$s = $ws.Range('D19').Style
$ws.Range('D19').Value = "'0.06559"
$ws.Range('D19').Style = $s
$ws.Range('E19').Value = '  -2.08%  '
$s = $ws.Range('D20').Style
$ws.Range('D20').Value = "'17.42"
$ws.Range('D20').Style = $s
$ws.Range('E20').Value = '  -3.75%  '
$s = $ws.Range('D21').Style
$ws.Range('D21').Value = "'1.001"
$ws.Range('D21').Style = $s
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = '29.001.30'
$ws.Range('E22').Value = '  -2.73%  '
$s = $ws.Range('D23').Style
$ws.Range('D23').Value = "'5.432"
$ws.Range('D23').Style = $s
$ws.Range('E23').Value = '  -3.57%  '
$s = $ws.Range('D24').Style
$ws.Range('D24').Value = "'11.45"
$ws.Range('D24').Style = $s
$ws.Range('E24').Value = '  +0.95%  '
$s = $ws.Range('D25').Style
$ws.Range('D25').Value = "'2.200"
$ws.Range('D25').Style = $s
$ws.Range('E25').Value = '  -3.63%  '
$ws.Range('D26').Value = '2.137.95'
$ws.Range('E26').Value = '  -1.64%  '
$s = $ws.Range('D27').Style
$ws.Range('D27').Value = "'156.25"
$ws.Range('D27').Style = $s
$ws.Range('E27').Value = '  -3.54%  '
$s = $ws.Range('D28').Style
$ws.Range('D28').Value = "'19.60"
$ws.Range('D28').Style = $s
$ws.Range('E28').Value = '  -3.16%  '
$s = $ws.Range('D29').Style
$ws.Range('D29').Value = "'2.089"
$ws.Range('D29').Style = $s
$ws.Range('E29').Value = '  -5.21%  '
$s = $ws.Range('D30').Style
$ws.Range('D30').Value = "'5.497"
$ws.Range('D30').Style = $s
$ws.Range('E30').Value = '  -2.85%  '
$s = $ws.Range('D31').Style
$ws.Range('D31').Value = "'117.48"
$ws.Range('D31').Style = $s
$ws.Range('E31').Value = '  -4.36%  '
$s = $ws.Range('D32').Style
$ws.Range('D32').Value = "'1.043"
$ws.Range('D32').Style = $s
$ws.Range('E32').Value = '  +1.31%  '
$s = $ws.Range('D33').Style
$ws.Range('D33').Value = "'0.09328"
$ws.Range('D33').Style = $s
$ws.Range('E33').Value = '  -3.47%  '
$s = $ws.Range('D34').Style
$ws.Range('D34').Value = "'1.409"
$ws.Range('D34').Style = $s
$ws.Range('E34').Value = '  -5.22%  '
$s = $ws.Range('D35').Style
$ws.Range('D35').Value = "'3.528"
$ws.Range('D35').Style = $s
$ws.Range('E35').Value = '  -4.19%  '
$s = $ws.Range('D36').Style
$ws.Range('D36').Value = "'5.301"
$ws.Range('D36').Style = $s
$ws.Range('E36').Value = '  -3.70%  '
$s = $ws.Range('D37').Style
$ws.Range('D37').Value = "'0.06059"
$ws.Range('D37').Style = $s
$ws.Range('E37').Value = '  -3.21%  '
$s = $ws.Range('D38').Style
$ws.Range('D38').Value = "'0.02229"
$ws.Range('D38').Style = $s
$ws.Range('E38').Value = '  -4.28%  '
$s = $ws.Range('D39').Style
$ws.Range('D39').Value = "'8.364"
$ws.Range('D39').Style = $s
$ws.Range('E39').Value = '  -4.68%  '
$s = $ws.Range('D40').Style
$ws.Range('D40').Value = "'1.172"
$ws.Range('D40').Style = $s
$ws.Range('E40').Value = '  -2.26%  '
$s = $ws.Range('D41').Style
$ws.Range('D41').Value = "'0.9994"
$ws.Range('D41').Style = $s
$ws.Range('E41').Value = '  -0.03%  '
$s = $ws.Range('D42').Style
$ws.Range('D42').Value = "'0.5795"
$ws.Range('D42').Style = $s
$ws.Range('E42').Value = '  -5.40%  '
$s = $ws.Range('D43').Style
$ws.Range('D43').Value = "'0.1826"
$ws.Range('D43').Style = $s
$ws.Range('E43').Value = '  -4.47%  '
$s = $ws.Range('D44').Style
$ws.Range('D44').Value = "'10.10"
$ws.Range('D44').Style = $s
$ws.Range('E44').Value = '  -6.29%  '
$s = $ws.Range('D45').Style
$ws.Range('D45').Value = "'1.258"
$ws.Range('D45').Style = $s
$ws.Range('E45').Value = '  -2.66%  '
$s = $ws.Range('D46').Style
$ws.Range('D46').Value = "'0.07522"
$ws.Range('D46').Style = $s
$ws.Range('E46').Value = '  +0.59%  '
$s = $ws.Range('D47').Style
$ws.Range('D47').Value = "'2.288"
$ws.Range('D47').Style = $s
$ws.Range('E47').Value = '  -2.63%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$s = $ws.Range('D48').Style
$ws.Range('D48').Value = "'12.02"
$ws.Range('D48').Style = $s
$ws.Range('E48').Value = '  -4.51%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$s = $ws.Range('D49').Style
$ws.Range('D49').Value = "'0.5460"
$ws.Range('D49').Style = $s
$ws.Range('E49').Value = '  -4.65%  '
$ws.Range('B50').Value = 'PaxosStandard'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$s = $ws.Range('D50').Style
$ws.Range('D50').Value = "'1.125"
$ws.Range('D50').Style = $s
$ws.Range('E50').Value = '  +12.55%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$s = $ws.Range('D51').Style
$ws.Range('D51').Value = "'1.904"
$ws.Range('D51').Style = $s
$ws.Range('E51').Value = '  -5.05%  '
